# "Working on the MSR mining."
#
# The "Strengths" row of the comparison table had the sentence
# "Discovering unforeseen applications of Trajectory." split across two
# runs ("Discovering unforeseen applications " + "of Trajectory."). Merge
# them back into the single run that carries the first run's formatting,
# and drop the now-redundant second run.

$p = $ppt.ActivePresentation

# --- (cosmetic) notes page size tweak recorded alongside the text fix ---
# PowerPoint keeps this on the Page Setup dialog; harmless to (try to) set
# even if this host treats the notes page geometry as fixed.
try {
    $ps = $p.PageSetup
    $ps.NotesPageWidth = 540
    $ps.NotesPageHeight = 732
} catch {
}

# --- the actual content edit -------------------------------------------
$slide = $p.Slides.Item(1)
$tbl = $slide.Shapes.Item(1).Table

# Row 5 / Column 2 (1-indexed) = the "Strengths" row, "Pilot study" column.
$cell = $tbl.Cell(5, 2)
$tf = $cell.Shape.TextFrame
$textRange = $tf.TextRange

# Paragraph 2 of that cell holds:
#   run1: "Discovering unforeseen applications "
#   run2: "of Trajectory."
$para = $textRange.Paragraphs(2, 1)

$mergedText = "Discovering unforeseen applications of Trajectory."
$firstRunLen = "Discovering unforeseen applications ".Length
$secondRunLen = "of Trajectory.".Length

$firstRun = $para.Characters(1, $firstRunLen)
$secondRun = $para.Characters($firstRunLen + 1, $secondRunLen)

# Put the full sentence into the first run (keeps its rPr: sz=1100, dirty=0,
# smtClean=0) and empty out the second run so it disappears entirely.
$firstRun.Text = $mergedText
$secondRun.Text = ""
